$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.123.45"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "2.636.68"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'596.54"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").Value = "'154.81"
$ws.Range("E6").Value = "  +0.80%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.544"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "2.635.68"
$ws.Range("E9").Value = "  +0.59%  "
$ws.Range("D10").Value = "'0.144"
$ws.Range("E10").Value = "  +7.53%  "
$ws.Range("E11").Value = "  -0.90%  "
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("E13").Value = "  +1.12%  "
$ws.Range("D14").Value = "'27.88"
$ws.Range("E14").Value = "  +1.11%  "
$ws.Range("D15").Value = "'0.0000191"
$ws.Range("E15").Value = "  +2.10%  "
$ws.Range("D16").Value = "3.116.67"
$ws.Range("D17").Value = "68.065.21"
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").Value = "2.648.64"
$ws.Range("E18").Value = "  +0.61%  "
$ws.Range("D19").Value = "'11.35"
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("D20").Value = "'362.80"
$ws.Range("E20").Value = "  -1.23%  "
$ws.Range("D21").Value = "'7.42"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("E22").Value = "  +3.05%  "
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("E24").Value = "  -1.02%  "
$ws.Range("D25").Value = "'74.85"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").Value = "'9.67"
$ws.Range("E27").Value = "  -2.70%  "
$ws.Range("E28").Value = "  +1.14%  "
$ws.Range("D29").Value = "2.769.50"
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("D31").Value = "'558.87"
$ws.Range("E31").Value = "  -2.14%  "
$ws.Range("E32").Value = "  +1.37%  "
$ws.Range("E33").Value = "  +0.42%  "
$ws.Range("E34").Value = "  +0.96%  "
$ws.Range("E35").Value = "  +1.63%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("D37").Value = "'1.56"
$ws.Range("E37").Value = "  +3.05%  "
$ws.Range("D38").Value = "'161.22"
$ws.Range("E38").Value = "  -0.73%  "
$ws.Range("D39").Value = "'19.29"
$ws.Range("E39").Value = "  +1.01%  "
$ws.Range("D40").Value = "'0.371"
$ws.Range("E40").Value = "  +1.29%  "
$ws.Range("E41").Value = "  -0.78%  "
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("D43").Value = "0.0₆0340"
$ws.Range("E43").Value = "  +5.33%  "
$ws.Range("E44").Value = "  +1.16%  "
$ws.Range("D45").Value = "'2.62"
$ws.Range("E45").Value = "  -0.91%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("E47").Value = "  +0.72%  "
$ws.Range("D48").Value = "'158.32"
$ws.Range("E48").Value = "  +2.30%  "
$ws.Range("E49").Value = "  +1.50%  "
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("E51").Value = "  +1.42%  "
